# Update workbook from version 0.9.15 (2025-10-30) to version 0.9.20 (2025-11-18)
# as described by the commit "feat: update to 0.9.20".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump Version and Date values.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.9.20"
$meta.Range("B8").Value = "2025-11-18T19:57:13-03:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: Schedule.actor no longer slices into
#    healthcareServices / locations. Instead it directly references the two
#    MedMe profiles, and a new invariant requires at least one Location.
# ---------------------------------------------------------------------------
$el = $wb.Worksheets.Item("Elements")

# Schedule.actor (row 34): update allowed reference types.
$el.Range("K34").Value = "Reference(https://fhir.medmehealth.com/pharmacy-services/StructureDefinition/medme-pharmacy-healthcareservice|https://fhir.medmehealth.com/pharmacy-services/StructureDefinition/medme-pharmacy-location)`n"

# Remove the slicing discriminator / description / rules - actor is no
# longer sliced.
$el.Range("AB34").ClearContents()
$el.Range("AC34").ClearContents()
$el.Range("AE34").ClearContents()

# Update the constraint(s) to add the new "at-least-one-location" invariant.
$el.Range("AJ34").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`nat-least-one-location:At least one actor must be a Location {where(reference.startsWith('Location/')).count() >= 1}"

# Remove the two actor slice rows (Schedule.actor:healthcareServices and
# Schedule.actor:locations). This shifts the following rows
# (Schedule.planningHorizon, Schedule.comment) up by two.
$el.Rows.Item(35).Delete()
$el.Rows.Item(35).Delete()

Write-Host "Applied 0.9.20 update"
